$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy formatting from the neighboring header cell (H1)
# so the new headers pick up the existing bold/centered/bordered style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: I and J values per row
$data = @(
    @(8, 8),   # row 2
    @(6, 6),   # row 3
    @(9, 9),   # row 4
    @(7, 9),   # row 5
    @(7, 7),   # row 6
    @(7, 7),   # row 7
    @(7, 7),   # row 8
    @(7, 7),   # row 9
    @(8, 8),   # row 10
    @(6, 7),   # row 11
    @(6, 6),   # row 12
    @(9, 9),   # row 13
    @(4, 5),   # row 14
    @(2, 2),   # row 15
    @(3, 4),   # row 16
    @(6, 6),   # row 17
    @(4, 4)    # row 18
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
